# Implemented Time Trend Analysis KPIs with engagement trends and MoM growth
#
# This script mirrors the removal of the "handle" column step that was
# recorded in the documentation workbook:
#  - "Gathering Info of Columns": the "handle" row is removed from the
#    after-transformation column list (G:H), shifting the remaining
#    entries up by one row.
#  - "Data Cleaning": a new "Remove Column" / "handle" step (row 12) is
#    logged, adding row 17 (+ a trailing blank row 18).
#  - "Data Wrangling": a matching "Remove Column" / "handle" step (row 9)
#    is logged as row 13.
#  - View/selection state is updated so that "Gathering Info of Columns"
#    becomes the active tab instead of "Story Telling - KPI's & Charts".

$wb = $excel.ActiveWorkbook

$wsInfo      = $wb.Worksheets.Item(1)   # Gathering Info of Columns
$wsCleaning  = $wb.Worksheets.Item(2)   # Data Cleaning
$wsWrangling = $wb.Worksheets.Item(3)   # Data Wrangling
$wsStory     = $wb.Worksheets.Item(4)   # Story Telling - KPI's & Charts

# ---------------------------------------------------------------------
# 1) "Gathering Info of Columns" - drop the "handle" row from the
#    Sr.No/Column Name summary table (columns G:H, rows 9-25) and shift
#    the remaining rows up by one (values + formatting).
# ---------------------------------------------------------------------
$wsInfo.Range("H10:H25").Copy($wsInfo.Range("H9:H24"))
$wsInfo.Range("H25").Clear()
$wsInfo.Range("G25").ClearContents()

# ---------------------------------------------------------------------
# 2) "Data Cleaning" - append the "Remove Column" / "handle" step as a
#    new row 17 (modeled on row 16), plus a trailing blank row 18.
# ---------------------------------------------------------------------
$wsCleaning.Range("B16:L16").Copy($wsCleaning.Range("B17:L17"))
$wsCleaning.Range("E17").Clear()

$wsCleaning.Cells.Item(17, 2).Value2 = 12
$wsCleaning.Cells.Item(17, 3).Value2 = "Remove Column"
$wsCleaning.Cells.Item(17, 4).Value2 = "handle"
$wsCleaning.Cells.Item(17, 6).Value2 = 12
$wsCleaning.Cells.Item(17, 7).Formula = "=H16"
$wsCleaning.Cells.Item(17, 8).Value2 = 90343
$wsCleaning.Cells.Item(17, 9).Formula = "=G17-H17"
$wsCleaning.Cells.Item(17, 10).Formula = "=K16"
$wsCleaning.Cells.Item(17, 11).Value2 = 19
$wsCleaning.Cells.Item(17, 12).Value2 = 1
$wsCleaning.Rows.Item(17).RowHeight = 15.6

$wsCleaning.Cells.Item(16, 2).Copy($wsCleaning.Cells.Item(18, 2))
$wsCleaning.Cells.Item(18, 2).ClearContents()
$wsCleaning.Rows.Item(18).RowHeight = 15.6

# ---------------------------------------------------------------------
# 3) "Data Wrangling" - append the matching "Remove Column" / "handle"
#    step as a new row 13 (modeled on row 5).
# ---------------------------------------------------------------------
$wsWrangling.Range("D5:F5").Copy($wsWrangling.Range("D13:F13"))
$wsWrangling.Cells.Item(13, 4).Value2 = 9
$wsWrangling.Cells.Item(13, 5).Value2 = "handle"
$wsWrangling.Cells.Item(13, 6).Value2 = "Remove Column"
$wsWrangling.Rows.Item(13).RowHeight = 15.6

# ---------------------------------------------------------------------
# 4) View state - update selections on the non-active sheets first,
#    then finish on "Gathering Info of Columns" so it becomes the
#    active/selected tab (moving away from "Story Telling").
# ---------------------------------------------------------------------
$wsCleaning.Range("D17").Select()
$wsWrangling.Range("E13").Select()

$wsInfo.Activate()
$wsInfo.Range("J26").Select()
